$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "missing-number" question row (row 7) ---
$ws.Range("A7").Value = "missing-number"
$ws.Range("C7").Value = "Use Cycle Sort"

# Add the hyperlink for the new LeetCode question (matches the pattern used
# by the other question cells: A2, A3, A4, A6). Passing the URL as the
# TextToDisplay keeps the <hyperlink display="..."> attribute equal to the
# address (as Excel does for the existing rows), and we restore the cell's
# visible text right after since Hyperlinks.Add overwrites it.
$ws.Hyperlinks.Add($ws.Range("A7"), "https://leetcode.com/problems/missing-number", "", "", "https://leetcode.com/problems/missing-number") | Out-Null
$ws.Range("A7").Value = "missing-number"

# Apply the same "Hyperlink" style used by the other linked question cells.
# Doing this AFTER the value/hyperlink are set lets the engine reuse the
# existing Hyperlink cell-format record instead of minting a new one.
$ws.Range("A7").Style = "Hyperlink"

# Move the active selection to the next empty row, matching the workbook
# state after the new row was entered.
$ws.Range("A8").Select()
